# Actualización automática del inventario, Google Sheets y productos.json
# Adds a new inventory row (row 74) to the single worksheet, matching the
# structure of the existing data rows (columns A,B,D,E,F,G,H,I,J - no C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

$ws.Cells.Item($row, 1).Value = "QHYH50"
$ws.Cells.Item($row, 2).Value = "Guantes conluz led"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 9
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E74-D74)*G74"
$ws.Cells.Item($row, 9).Formula = "=D74*F74"
$ws.Cells.Item($row, 10).Value = 0
